$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 13 de Octubre de 2020 a las 22:22"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 8073700
$ws.Range("C4").Value = 35374
$ws.Range("D4").Value = 5208928
$ws.Range("E4").Value = 2644221
$ws.Range("G4").Value = 533
$ws.Range("H4").Value = 220551

# Sudafrica (row 14)
$ws.Range("B14").Value = 694537
$ws.Range("C14").Value = 1178
$ws.Range("D14").Value = 625574
$ws.Range("E14").Value = 50935
$ws.Range("G14").Value = 165
$ws.Range("H14").Value = 18028

# Alemania (row 25)
$ws.Range("B25").Value = 335527
$ws.Range("C25").Value = 4433
$ws.Range("E25").Value = 46687
$ws.Range("G25").Value = 19
$ws.Range("H25").Value = 9740

# Canada (row 30)
$ws.Range("B30").Value = 185371
$ws.Range("C30").Value = 2532
$ws.Range("D30").Value = 156447
$ws.Range("E30").Value = 19279
$ws.Range("G30").Value = 18
$ws.Range("H30").Value = 9645

# Costa de Marfil (row 92)
$ws.Range("B92").Value = 20183
$ws.Range("C92").Value = 28
$ws.Range("D92").Value = 19860
$ws.Range("E92").Value = 203

# Maldivas (row 105)
$ws.Range("B105").Value = 10993
$ws.Range("C105").Value = 50
$ws.Range("D105").Value = 9833

# Rows 108/109: Mozambique overtakes Guayana Francesa in total cases, swapping rank order.
# Row 108 now holds Mozambique's (updated) figures.
$ws.Range("A108").Value = "Mozambique"
$ws.Range("B108").Value = 10258
$ws.Range("C108").Value = 170
$ws.Range("D108").Value = 7880
$ws.Range("E108").Value = 2305
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 73

# Row 109 now holds Guayana Francesa's (unchanged) figures.
$ws.Range("A109").Value = "Guayana Francesa"
$ws.Range("B109").Value = 10192
$ws.Range("C109").Value = 12
$ws.Range("D109").Value = 9881
$ws.Range("E109").Value = 242
$ws.Range("H109").Value = 69

# row 118
$ws.Range("B118").Value = 7254
$ws.Range("C118").Value = 99
$ws.Range("D118").Value = 6210
$ws.Range("E118").Value = 967
$ws.Range("G118").Value = 2
$ws.Range("H118").Value = 77

# row 119
$ws.Range("B119").Value = 6680
$ws.Range("C119").Value = 192
$ws.Range("D119").Value = 2761
$ws.Range("E119").Value = 3697
$ws.Range("G119").Value = 3
$ws.Range("H119").Value = 222

# row 124
$ws.Range("B124").Value = 5696
$ws.Range("C124").Value = 13
$ws.Range("D124").Value = 5337
$ws.Range("E124").Value = 245
$ws.Range("G124").Value = 1
$ws.Range("H124").Value = 114

# row 125
$ws.Range("B125").Value = 5428
$ws.Range("C125").Value = 2
$ws.Range("D125").Value = 5362

# row 134
$ws.Range("B134").Value = 4908
$ws.Range("C134").Value = 3
$ws.Range("D134").Value = 4130
$ws.Range("E134").Value = 746

# row 136
$ws.Range("B136").Value = 4826
$ws.Range("C136").Value = 52
$ws.Range("D136").Value = 1364
$ws.Range("E136").Value = 3231
$ws.Range("G136").Value = 3
$ws.Range("H136").Value = 231

# row 138
$ws.Range("B138").Value = 4229
$ws.Range("C138").Value = 32
$ws.Range("D138").Value = 3849
$ws.Range("E138").Value = 347

# row 157
$ws.Range("B157").Value = 2309
$ws.Range("C157").Value = 3
$ws.Range("D157").Value = 1738
$ws.Range("E157").Value = 499

# row 166
$ws.Range("B166").Value = 1321
$ws.Range("C166").Value = 13
$ws.Range("D166").Value = 1120
$ws.Range("E166").Value = 109

# row 180
$ws.Range("B180").Value = 496
$ws.Range("C180").Value = 1
$ws.Range("D180").Value = 484
$ws.Range("E180").Value = 5

# row 205
$ws.Range("D205").Value = 30
$ws.Range("E205").Value = 0
